$d = $word.ActiveDocument

$replacements = @(
    @("2025-08-07 Thursday", "2025-08-08 Friday"),
    @("60×37=2220", "61×48=2928"),
    @("43×80=3440", "11×95=1045"),
    @("48×35=1680", "69×81=5589"),
    @("87×66=5742", "92×84=7728"),
    @("30×61=1830", "33×85=2805"),
    @("30×19=570", "88×71=6248"),
    @("95×77=7315", "95×48=4560"),
    @("14×58=812", "44×98=4312"),
    @("19×69=1311", "61×49=2989"),
    @("47×29=1363", "86×18=1548"),
    @("96×85=8160", "29×79=2291"),
    @("99×26=2574", "25×44=1100"),
    @("17×26=442", "38×48=1824"),
    @("42×16=672", "63×75=4725"),
    @("77×60=4620", "39×16=624"),
    @("87×68=5916", "49×97=4753"),
    @("53×32=1696", "61×38=2318"),
    @("93×28=2604", "62×60=3720"),
    @("66×55=3630", "30×25=750"),
    @("66×71=4686", "57×40=2280"),
    @("95×66=6270", "54×73=3942"),
    @("53×97=5141", "93×37=3441"),
    @("65×69=4485", "98×38=3724"),
    @("23×57=1311", "60×84=5040"),
    @("81×75=6075", "41×97=3977")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $range = $d.Content
    $range.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}
